# Insert two new columns (M_TotalTax, M_CorpTax) right after the M_POP
# column (column E) and before the GFA/IMF/OECD columns (old F:M, now
# shifted to H:O). Then populate the header + data for the new columns,
# and fix up the one data value that changed in column E (M_POP) for the
# UMICs row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing columns F:M (GFA/IMF/OECD data) two columns to the right
# by inserting two blank columns at F:G.
$ws.Range("F:G").Insert()

# New header cells
$ws.Range("F1").Value = "M_TotalTax"
$ws.Range("G1").Value = "M_CorpTax"

# New column data (M_TotalTax, M_CorpTax) for each group row
$ws.Range("F2").Value = 14106286460237.92
$ws.Range("G2").Value = 1155021202746.413

$ws.Range("F3").Value = 3207987015.574299
$ws.Range("G3").Value = 0

$ws.Range("F4").Value = 734615892234.8064
$ws.Range("G4").Value = 88889835996.30263

$ws.Range("F5").Value = 558865056646.082
$ws.Range("G5").Value = 72600947639.16805

$ws.Range("F6").Value = 4579473077980.816
$ws.Range("G6").Value = 674619880691.7614

# Corrected M_POP value for the UMICs row (column E, row 6)
$ws.Range("E6").Value = 2427884184.75
